# Update metrics table (columns B..Q, rows 2..26) with new values from the
# re-trained model ("atualizado todo o treinamento para o novo lm").
# Every data row shares the same new metric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B" = 0.9999968104599276
    "C" = 0.9990493596125468
    "D" = 0.9999999759600751
    "E" = 0.9999801781834854
    "F" = 0.9999961222740886
    "G" = [double]"2.97729565298149e-06"
    "H" = 0.0008873810733972544
    "I" = [double]"2.938460367545546e-08"
    "J" = [double]"7.012419020754245e-06"
    "K" = [double]"3.52090181221485e-06"
    "L" = 0.0001188178455764825
    "M" = 0.001725484179290407
    "N" = 1.000004502880102
    "O" = 0.001798941613047459
    "P" = 107.4489903369459
    "Q" = 157.4228991565422
}

for ($row = 2; $row -le 26; $row++) {
    foreach ($col in $newValues.Keys) {
        $ws.Range("$col$row").Value = $newValues[$col]
    }
}
